$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full contents of row 7 and row 8 (all populated
# columns, A:AI) - every field that differed between the two rows is
# exchanged; columns that already held identical values in both rows
# (C, K, T, U, V, W, AD, AE, AG, AH, AT, AW, AX, AY) are left untouched.

# --- Plain numeric columns: straight value swap -----------------------
$numericCols = 1, 2, 5, 17, 18, 19   # A, B, E, Q, R, S

foreach ($col in $numericCols) {
    $c7 = $ws.Cells.Item(7, $col)
    $c8 = $ws.Cells.Item(8, $col)
    $v7 = $c7.Value2
    $v8 = $c8.Value2
    $c7.Value2 = $v8
    $c8.Value2 = $v7
}

# --- Plain text columns (not numeric- or date-looking): swap directly -
$textCols = 4, 6, 7, 8, 16   # D, F, G, H, P

foreach ($col in $textCols) {
    $c7 = $ws.Cells.Item(7, $col)
    $c8 = $ws.Cells.Item(8, $col)
    $v7 = $c7.Value2
    $v8 = $c8.Value2
    $c7.Value2 = $v8
    $c8.Value2 = $v7
}

# --- Text columns whose content looks like a date/time/number: force
#     a text number-format first so Excel doesn't reinterpret the
#     swapped string as a date serial / number -------------------------
$textGuardCols = 25, 26, 27, 28   # Y, Z, AA, AB (Startdatum/tid, Slutdatum/tid)

foreach ($col in $textGuardCols) {
    $c7 = $ws.Cells.Item(7, $col)
    $c8 = $ws.Cells.Item(8, $col)
    $v7 = $c7.Value2
    $v8 = $c8.Value2
    $c7.NumberFormat = "@"
    $c8.NumberFormat = "@"
    $c7.Value2 = $v8
    $c8.Value2 = $v7
    # Restore the original (General) style now that the literal text is
    # safely stored, so no spurious style/number-format diff remains.
    $c7.Style = "Normal"
    $c8.Style = "Normal"
}

# --- I7/I8 ("Antal"): text "70" <-> blank -------------------------------
$i7 = $ws.Cells.Item(7, 9)
$i8 = $ws.Cells.Item(8, 9)
$i7Val = $i7.Value2
$i8.NumberFormat = "@"
$i8.Value2 = $i7Val
$i8.Style = "Normal"
$i7.ClearContents()

# --- AI7/AI8 ("Biotop-beskrivning"): "Kalkbarrskog" <-> blank -----------
$ai7 = $ws.Cells.Item(7, 35)
$ai8 = $ws.Cells.Item(8, 35)
$ai8Val = $ai8.Value2
$ai7.Value2 = $ai8Val
$ai8.ClearContents()

# --- J8/N8/AF8 ("Enhet"/"Metod"/"Bestämningsmetod"): become blank -------
# (J7/N7/AF7 are already blank and gain no visible content, so nothing
# further is required on the row-7 side.)
$ws.Cells.Item(8, 10).ClearContents()   # J8
$ws.Cells.Item(8, 14).ClearContents()   # N8
$ws.Cells.Item(8, 32).ClearContents()   # AF8
